# feat: add 2022-Q4 data
#
# - total sheet ("总计"): duplicate the existing 2022-Q1 row as a new row,
#   then overwrite the original row in-place with the new 2022-Q4 totals.
# - worksheets: clone the existing "2022-Q1" sheet (preserving its data byte
#   for byte) into a brand new sheet placed right after it, then repurpose the
#   original sheet object to hold the new "2022-Q4" fund data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "总计"
$ws2 = $wb.Worksheets.Item(2)   # "2022-Q1" (will be turned into "2022-Q4")

# Helper: write a value that must stay TEXT even though it looks numeric
# ("001420", "1.71", ...) without leaving a stray NumberFormat-driven style
# behind on the cell.
function Set-TextValue {
    param($Cell, $Val)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Val
    $Cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. "总计" sheet: push the current 2022-Q1 row down to row 3, keeping its
#    formatting/value intact, and give it the new index (1).
# ---------------------------------------------------------------------
$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))
$ws1.Cells.Item(3, 1).Value = 1

# ---------------------------------------------------------------------
# 2. Clone the "2022-Q1" worksheet into a new sheet right after it so the
#    original quarterly fund-holding data survives unchanged.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws2)
$ws2.UsedRange.Copy($newSheet.Range("A1"))

$ws2.Name = "2022-Q4"
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 3. Now fill in row 2 of "总计" with the new 2022-Q4 summary values
#    (row index / col A stays 0, as it was).
# ---------------------------------------------------------------------
$ws1.Cells.Item(2, 2).Value = "2022-Q4"
$ws1.Cells.Item(2, 3).Value = 2
$ws1.Cells.Item(2, 4).Value = 0.01

# ---------------------------------------------------------------------
# 4. Replace the (now renamed) "2022-Q4" sheet's contents with the new
#    fund-holding data.
# ---------------------------------------------------------------------
$ws2.Cells.Clear()

# Stamp the header row + index column with the same style used for the
# index column on the "总计" sheet (cellXfs index 2: centered/top, bordered).
$ws1.Range("A2").Copy($ws2.Range("B1"))
$ws1.Range("A2").Copy($ws2.Range("C1"))
$ws1.Range("A2").Copy($ws2.Range("D1"))
$ws1.Range("A2").Copy($ws2.Range("E1"))
$ws1.Range("A2").Copy($ws2.Range("F1"))
$ws1.Range("A2").Copy($ws2.Range("G1"))
$ws1.Range("A2").Copy($ws2.Range("H1"))
$ws1.Range("A2").Copy($ws2.Range("A2"))
$ws1.Range("A2").Copy($ws2.Range("A3"))

$ws2.Cells.Item(1, 2).Value = "基金代码"
$ws2.Cells.Item(1, 3).Value = "基金名称"
$ws2.Cells.Item(1, 4).Value = "基金规模"
$ws2.Cells.Item(1, 5).Value = "股票总仓位"
$ws2.Cells.Item(1, 6).Value = "仓位占比"
$ws2.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2.Cells.Item(1, 8).Value = "仓位排名"

$ws2.Cells.Item(2, 1).Value = 0
Set-TextValue $ws2.Cells.Item(2, 2) "001420"
$ws2.Cells.Item(2, 3).Value = "南方大数据300指数A"
Set-TextValue $ws2.Cells.Item(2, 4) "1.71"
Set-TextValue $ws2.Cells.Item(2, 5) "93.56"
Set-TextValue $ws2.Cells.Item(2, 6) "0.67"
Set-TextValue $ws2.Cells.Item(2, 7) "0.0115"
$ws2.Cells.Item(2, 8).Value = 10

$ws2.Cells.Item(3, 1).Value = 1
Set-TextValue $ws2.Cells.Item(3, 2) "001426"
$ws2.Cells.Item(3, 3).Value = "南方大数据300指数C"
Set-TextValue $ws2.Cells.Item(3, 4) "0.32"
Set-TextValue $ws2.Cells.Item(3, 5) "93.56"
Set-TextValue $ws2.Cells.Item(3, 6) "0.67"
Set-TextValue $ws2.Cells.Item(3, 7) "0.0021"
$ws2.Cells.Item(3, 8).Value = 10
